$wb = $excel.ActiveWorkbook

# --- unit_file: add two new manual-correction rows for plant 50489 ---
# (C4 and C5 generators get a new "prop" = 0 correction, in addition to
# the existing "prime_mover" = GT corrections already present)
$ws2 = $wb.Worksheets.Item("unit_file")

$ws2.Cells.Item(24, 1).Value = "50489"
$ws2.Cells.Item(24, 2).Value = "C4"
$ws2.Cells.Item(24, 4).Value = "prop"
$ws2.Cells.Item(24, 5).Value = "0"

$ws2.Cells.Item(25, 1).Value = "50489"
$ws2.Cells.Item(25, 2).Value = "C5"
$ws2.Cells.Item(25, 4).Value = "prop"
$ws2.Cells.Item(25, 5).Value = "0"

# --- plant_file: remove the old primary_fuel_type / primary_fuel_category
# manual corrections for plants 55970 and 10154 (rows 10-13) ---
$ws3 = $wb.Worksheets.Item("plant_file")
$ws3.Rows("10:13").Delete()

# --- update selection / active sheet state to match the final workbook ---
$ws2.Activate()
$ws2.Range("E25").Select()

$ws3.Activate()
$ws3.Range("C14").Select()
